$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230, shifting existing row 230 (and below) down to 231.
$ws.Rows.Item(230).Insert()

# Populate the new row 230 with values. Columns A,B,C,E,F,G,H,I,N,O,Q,R are the same
# as the row that used to be at 230 (now at 231); only D, J, K, L, M, P differ.
$ws.Cells.Item(230, 1).Value = 4
$ws.Cells.Item(230, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(230, 3).Value = "Los Lagos"
$ws.Cells.Item(230, 4).Value = 44468
$ws.Cells.Item(230, 5).Value = 10
$ws.Cells.Item(230, 6).Value = 100112033
$ws.Cells.Item(230, 7).Value = "Lechuga"
$ws.Cells.Item(230, 8).Value = "Escarola"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 150
$ws.Cells.Item(230, 11).Value = 12000
$ws.Cells.Item(230, 12).Value = 12000
$ws.Cells.Item(230, 13).Value = 12000
$ws.Cells.Item(230, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(230, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(230, 16).Value = 800
$ws.Cells.Item(230, 17).Value = 15
$ws.Cells.Item(230, 18).Value = "Hortaliza"
